$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.712.06'
$ws.Range("E2").Value = '  -1.93%  '
$ws.Range("D3").Value = '1.806.84'
$ws.Range("E3").Value = '  -2.22%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'231.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").Value = "'0.604"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = "'39.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -6.53%  '
$ws.Range("D9").Value = "'0.318"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.99%  '
$ws.Range("D10").Value = "'0.0680"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.61%  '
$ws.Range("D11").Value = "'0.0991"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.80%  '
$ws.Range("D12").Value = '2.066.11'
$ws.Range("E12").Value = '  -2.35%  '
$ws.Range("D13").Value = '1.792.82'
$ws.Range("E13").Value = '  -3.30%  '
$ws.Range("D14").Value = "'0.661"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.91%  '
$ws.Range("D15").Value = "'10.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.98%  '
$ws.Range("D16").Value = "'4.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.35%  '
$ws.Range("D17").Value = '34.696.81'
$ws.Range("E17").Value = '  -2.01%  '
$ws.Range("D18").Value = "'69.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("D19").Value = '0.0₃0783'
$ws.Range("E19").Value = '  -1.85%  '
$ws.Range("D20").Value = "'239.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.96%  '
$ws.Range("D21").Value = "'11.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.04%  '
$ws.Range("D22").Value = "'4.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.46%  '
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("E24").Value = '  +2.19%  '
$ws.Range("D25").Value = "'172.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.07%  '
$ws.Range("D26").Value = "'7.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.31%  '
$ws.Range("D27").Value = "'17.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.00%  '
$ws.Range("E28").Value = '  -2.08%  '
$ws.Range("E29").Value = '  +9.36%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = "'3.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.87%  '
$ws.Range("D32").Value = "'0.0545"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("D33").Value = "'3.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.26%  '
$ws.Range("D34").Value = "'1.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +15.40%  '
$ws.Range("E35").Value = '  -5.13%  '
$ws.Range("D36").Value = "'0.695"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.99%  '
$ws.Range("D37").Value = "'91.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.18%  '
$ws.Range("D38").Value = "'1.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.28%  '
$ws.Range("D39").Value = '1.310.49'
$ws.Range("E39").Value = '  -3.07%  '
$ws.Range("D40").Value = "'0.0191"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.71%  '
$ws.Range("E41").Value = '  -0.90%  '
$ws.Range("D42").Value = "'0.959"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.87%  '
$ws.Range("D43").Value = "'14.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.76%  '
$ws.Range("D44").Value = "'2.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -9.11%  '
$ws.Range("D45").Value = "'2.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.50%  '
$ws.Range("D46").Value = "'6.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.59%  '
$ws.Range("D47").Value = "'0.0511"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.63%  '
$ws.Range("D48").Value = '1.994.08'
$ws.Range("E48").Value = '  -1.04%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = "'0.0671"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.08%  '
$ws.Range("B50").Value = 'PaxDollar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D50").Value = "'1.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("D51").Value = "'98.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.43%  '
